$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Childnodes" column (D) values for rows 33-36
$ws.Range("D33").Value = 25
$ws.Range("D34").Value = 25
$ws.Range("D35").Value = 25
$ws.Range("D36").Value = 27
